$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new line-item row (J1 - ARM Cortex-M 10pin connector)
$ws.Range("G16").Value = "https://www.digikey.ca/products/en?keywords=%201175-1629-ND"
$ws.Range("A16").Value = "ARM Cortex-M 10pin conn"
$ws.Range("B16").Value = "10pin connector"
$ws.Range("C16").Value = "J1"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 0.082
$ws.Range("H16").Formula = "=E16*F16"

# Update the view: scrolled to column C, selection moved to J20
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J20").Select()
